$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing row 195 (old weekly record) down into new row 197,
# and the existing row 196 (old weekly record) down into new row 198,
# preserving all of their original values before row 195/196 get overwritten
# with the new week's data.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(197, $col).Value = $ws.Cells.Item(195, $col).Value2
    $ws.Cells.Item(198, $col).Value = $ws.Cells.Item(196, $col).Value2
}

# Copy the date style (format) used in column D down to the new rows too.
$ws.Cells.Item(195, 4).Copy() | Out-Null
$ws.Cells.Item(197, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(196, 4).Copy() | Out-Null
$ws.Cells.Item(198, 4).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Now overwrite row 195 with the new week's values.
$ws.Cells.Item(195, 4).Value = 44595
$ws.Cells.Item(195, 10).Value = 40
$ws.Cells.Item(195, 11).Value = 7000
$ws.Cells.Item(195, 12).Value = 7000
$ws.Cells.Item(195, 13).Value = 7000
$ws.Cells.Item(195, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(195, 16).Value = 7000

# And overwrite row 196 with the new week's values.
$ws.Cells.Item(196, 4).Value = 44595
$ws.Cells.Item(196, 10).Value = 100
$ws.Cells.Item(196, 11).Value = 7000
$ws.Cells.Item(196, 12).Value = 7000
$ws.Cells.Item(196, 13).Value = 7000
$ws.Cells.Item(196, 16).Value = 7000
